# The deck's handout master has a "datetimeFigureOut" date field
# (Insert > Header & Footer > Date and time) currently cached as
# 10/31/2022. Update it to 11/17/2022 via the documented
# HeadersFooters object, the supported automation surface for
# editing header/footer placeholder text (including the date).
$p  = $ppt.ActivePresentation
$hm = $p.HandoutMaster
$hm.HeadersFooters.DateAndTime.Text = "11/17/2022"
